$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.035.18"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "1.964.98"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.90"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4967"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4208"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.03"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09215"
$ws.Range("E10").Value = "  +4.40%  "
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.80"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "2.004.94"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.881"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.457"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.77"
$ws.Range("E17").Value = "  -4.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001102"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.23"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.947"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "29.069.81"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.99"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "2.227.34"
$ws.Range("E26").Value = "  -6.60%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.34"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.330"
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.251"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.56"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.046"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09837"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.518"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.815"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.716"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.323"
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.025"
$ws.Range("E39").Value = "  -6.40%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06367"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6448"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1985"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6203"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.350"
$ws.Range("E46").Value = "  +6.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.194"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.470"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000324"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06980"
$ws.Range("E51").Value = "  -0.63%  "
